$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.624.03"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "1.626.79"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.487"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "1.855.86"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "1.624.01"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "26.611.37"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "206.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("D36").Value = "1.161.85"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.500"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.788"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "1.766.02"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
